$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.693.51'
$ws.Range('E2').Value = '  -0.16%  '
$ws.Range('D3').Value = '2.209.84'
$ws.Range('E3').Value = '  -2.40%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '230.11'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.75%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.613'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.98%  '
$ws.Range('E7').Value = '  -6.26%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.403'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.47%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '57.11'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -5.12%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0889'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.02%  '
$ws.Range('E12').Value = '  -2.57%  '
$ws.Range('D13').Value = '2.539.13'
$ws.Range('E13').Value = '  -2.25%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '15.46'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -4.37%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '22.22'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.56%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.67'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.65%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.795'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -4.12%  '
$ws.Range('D18').Value = '2.213.39'
$ws.Range('E18').Value = '  -2.29%  '
$ws.Range('D19').Value = '41.652.17'
$ws.Range('E19').Value = '  +0.06%  '
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').Value = '0.0₃0902'
$ws.Range('E20').Value = '  -4.26%  '
$ws.Range('B21').Value = 'Litecoin'
$ws.Range('C21').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '72.15'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.73%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.06'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.35%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '241.04'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.59%  '
$ws.Range('E25').Value = '  -2.83%  '
$ws.Range('E26').Value = '  -3.33%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.68'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.01%  '
$ws.Range('E28').Value = '  -1.10%  '
$ws.Range('E29').Value = '  -5.58%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.46'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.32%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '19.78'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.93%  '
$ws.Range('E32').Value = '  -8.59%  '
$ws.Range('E33').Value = '  -3.78%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.02'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.88%  '
$ws.Range('E35').Value = '  -3.26%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0650'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.46%  '
$ws.Range('E37').Value = '  -4.84%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.32'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -8.70%  '
$ws.Range('E39').Value = '  -8.55%  '
$ws.Range('E40').Value = '  -8.08%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.999'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.21%  '
$ws.Range('E42').Value = '  -1.19%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.57'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.19%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0957'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.42%  '
$ws.Range('B45').Value = 'TrustWalletToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.20'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.51%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '97.30'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.47%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.39'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -15.00%  '
$ws.Range('D48').Value = '1.468.50'
$ws.Range('E48').Value = '  -2.74%  '
$ws.Range('E49').Value = '  -7.34%  '
$ws.Range('E50').Value = '  -0.99%  '
$ws.Range('E51').Value = '  -5.59%  '
